$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header and description values in the specific order needed so that the
# shared-strings table ends up built in the same order as the source edit.
$ws.Range("F1").Value = "Description"
$ws.Range("F2").Value = "Create start window in application with login and password inputs"
$ws.Range("F3").Value = "Create main window with setting fields from Performanse Counters and buttons Start and Stop"
$ws.Range("F4").Value = "Create model which includes main application work"
$ws.Range("F5").Value = "Create subodules ""Managers"" from control every parts application"
$ws.Range("F6").Value = "Add initialization select Performance Counters"
$ws.Range("F7").Value = "Add cached files in application"
$ws.Range("F8").Value = "Add complex Performance Couters from work with Matworl Loader"
$ws.Range("F11").Value = "Add class from work with system processes"
$ws.Range("F9").Value = "Create logger from created log files"
$ws.Range("F10").Value = "Beatiful display Message Box in logger class"
$ws.Range("F12").Value = "Add window from realization function"
$ws.Range("F13").Value = "Add hidden main window in Windows start panel"
$ws.Range("F14").Value = "Add validation from input fields in main window"
$ws.Range("F15").Value = "Create console application from receive messages from the client"
$ws.Range("F16").Value = "Connect EnityFramework to application and create database model"
$ws.Range("F17").Value = "Add tools from registration user in database"
$ws.Range("F18").Value = "Add conditions on main window from separate user logic"
$ws.Range("F19").Value = "Add ResourceDictinary from plugins"
$ws.Range("F20").Value = "Create new window from check currentsystem data"
$ws.Range("F21").Value = "Add checkbox on main window from setting disabling message box"
$ws.Range("F22").Value = "Add resx files from localization app"
$ws.Range("F23").Value = "Create VisualStudio Installers from Server and Application"
$ws.Range("F24").Value = "Added upload current user setting to database"
$ws.Range("F25").Value = "Add new folder in AppData from each user. Create uniqle log file from him"
$ws.Range("F26").Value = "Create uniqle cach file from each user"
$ws.Range("F27").Value = "Create block from check USB device"
$ws.Range("F28").Value = "Add Name current user in windows Title"
$ws.Range("F29").Value = "Update validation from time fields to 5 seconds"
$ws.Range("F30").Value = "Add validation on input field from Process Windows"
$ws.Range("F31").Value = "Refactoring code on new branch"

$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("I33").Select()
